# Update data of httpd
# - Adds a new "Max length" column (J) with values on Sheet2
# - Adds a new "httpd" benchmark row (row 18) on Sheet2
# - Switches the active sheet/selection from Sheet3 back to Sheet2

$wb = $excel.ActiveWorkbook

$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

# --- New "Max length" column header ---
$ws2.Range("J1").Value = "Max length"

# --- Fill in "Max length" values for existing rows ---
$ws2.Range("J2").Value = 6
$ws2.Range("J3").Value = 1
$ws2.Range("J4").Value = 7
$ws2.Range("J5").Value = 5
$ws2.Range("J6").Value = 1
$ws2.Range("J7").Value = 1
$ws2.Range("J8").Value = 4
$ws2.Range("J9").Value = 8
$ws2.Range("J10").Value = 1
$ws2.Range("J11").Value = 7

# --- New "httpd" benchmark row ---
$ws2.Range("A18").Value = "httpd"
$ws2.Range("B18").Value = 398
$ws2.Range("C18").Value = 233
$ws2.Range("D18").Value = 371
$ws2.Range("E18").Value = 53
$ws2.Range("F18").Value = "{0: 40, 1: 13}"
$ws2.Range("G18").Value = "{0: 216, 1: 155}"
$ws2.Range("H18").Value = "{0: 220, 1: 13}"
$ws2.Range("I18").Value = 1
$ws2.Range("J18").Value = 6

# --- Match styling of neighbouring rows for the new row ---
# (A18/B18.../J18 keep the default/general style, like rows 15-16;
#  F18:H18 use the wrapped-text style shared by the F/G/H columns above)
$ws2.Range("F18:H18").Style = $ws2.Range("F8:H8").Style

# --- Switch active sheet / selection back to Sheet2 ---
$ws3.Range("B1").Select()
$ws2.Activate()
$ws2.Range("I19").Select()
